$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Previous closing" (B), "Closing" (C) and "Change (%)" (D) columns
# for each sector row (rows 2-12) with refreshed BRVM index figures.

$ws.Range("B2").Value = 109.69
$ws.Range("C2").Value = 110.12
$ws.Range("D2").Value = 0.39

$ws.Range("B3").Value = 158.53
$ws.Range("C3").Value = 157.81
$ws.Range("D3").Value = -0.45

$ws.Range("B4").Value = 937.75
$ws.Range("C4").Value = 869.89
$ws.Range("D4").Value = -7.24

$ws.Range("B5").Value = 219.4
$ws.Range("C5").Value = 220.13
$ws.Range("D5").Value = 0.33

$ws.Range("B6").Value = 322.08
$ws.Range("C6").Value = 319.87
$ws.Range("D6").Value = -0.69

$ws.Range("B7").Value = 87.13
$ws.Range("C7").Value = 88.41
$ws.Range("D7").Value = 1.47

$ws.Range("B8").Value = 104.7
$ws.Range("C8").Value = 104.46
$ws.Range("D8").Value = -0.23

$ws.Range("B9").Value = 102.24
$ws.Range("C9").Value = 101.99
$ws.Range("D9").Value = -0.24

$ws.Range("B10").Value = 106.58
$ws.Range("C10").Value = 107.45
$ws.Range("D10").Value = 0.82

$ws.Range("B11").Value = 554.22
$ws.Range("C11").Value = 553.17
$ws.Range("D11").Value = -0.19

$ws.Range("B12").Value = 348.86
$ws.Range("C12").Value = 342.9
$ws.Range("D12").Value = -1.71
